$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The report previously listed 3 facilities (rows 2-4). Two more facilities
# ("Crumpet" and "Scone") are being added as rows 5 and 6, re-using the same
# formatting as the existing facility rows above them.

# Copy the formatting (number formats / fonts / fill) of the last existing
# data row (row 4, columns A:J) down onto the two new rows so the new cells
# pick up the same styles used throughout the table.
$ws.Range("A4:J4").Copy() | Out-Null
$ws.Range("A5:J5").PasteSpecial(-4122) | Out-Null
$ws.Range("A6:J6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 5 - Crumpet GEF
$ws.Range("A5").Value = "Crumpet GEF"
$ws.Range("B5").Value = 20001371
$ws.Range("C5").Value = "Crumpet exporter"
$ws.Range("D5").Value = "GBP"
$ws.Range("E5").Value = 7000000
$ws.Range("F5").Value = 3938753.8
$ws.Range("G5").Value = 777
$ws.Range("H5").Value = 456
$ws.Range("I5").Value = "GBP"
$ws.Range("J5").Value = "GBP"

# Row 6 - Scone GEF
$ws.Range("A6").Value = "Scone GEF"
$ws.Range("B6").Value = 20001371
$ws.Range("C6").Value = "Scone exporter"
$ws.Range("D6").Value = "GBP"
$ws.Range("E6").Value = 770000
$ws.Range("F6").Value = 761579.37
$ws.Range("G6").Value = 777
$ws.Range("H6").Value = 456.77
$ws.Range("I6").Value = "GBP"
$ws.Range("J6").Value = "GBP"

# Reflect the user's final selection: the two newly entered rows.
$ws.Activate() | Out-Null
$ws.Range("A5:J6").Select() | Out-Null
